$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the formatting of the existing header row
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column data
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
